function Insert-RunsAsNewParagraphAfter {
    param($insertionPoint, $runsXml)
    $pkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $runsXml + '</w:p><w:sectPr><w:pgSz w:w="12240" w:h="15840"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $insertionPoint.InsertXML($pkg) | Out-Null
}

# Appends one or more <w:r>...</w:r> runs (raw WordML) to the end of $paragraph,
# as genuinely separate runs (not merged with the existing trailing run), by
# inserting them as a new paragraph directly after $paragraph and then
# deleting the paragraph mark that separates them, splicing the two
# paragraphs back into one.
function Append-Runs {
    param($doc, $paraIndex, $runsXml)
    $p = $doc.Paragraphs.Item($paraIndex)
    $ip = $doc.Range($p.Range.End - 1, $p.Range.End - 1)
    Insert-RunsAsNewParagraphAfter $ip $runsXml
    $p = $doc.Paragraphs.Item($paraIndex)
    $mark = $doc.Range($p.Range.End - 1, $p.Range.End)
    $mark.Delete()
}

$d = $word.ActiveDocument

# --- Paragraph 1: "SE16N"  ->  "SE38" + <tab/> + " Program oluşturma." ---
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Text = "SE38"
Append-Runs $d 1 '<w:r><w:tab/><w:t xml:space="preserve"> Program oluşturma.</w:t></w:r>'

# --- Paragraph 2: "SE38"  ->  " " + "SE11" + <tab/> + " Tablo,domain,data element" + ",structure" + "." ---
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Text = " "
Append-Runs $d 2 '<w:r><w:t>SE11</w:t></w:r>'
Append-Runs $d 2 '<w:r><w:tab/><w:t xml:space="preserve"> Tablo,domain,data element</w:t></w:r>'
Append-Runs $d 2 '<w:r><w:t>,structure</w:t></w:r>'
Append-Runs $d 2 '<w:r><w:t>.</w:t></w:r>'

# --- Paragraph 3: "SE11"  ->  "SE16N" + " Tablo görüntüleme,editleme" ---
$p3 = $d.Paragraphs.Item(3)
$p3.Range.Text = "SE16N"
Append-Runs $d 3 '<w:r><w:t xml:space="preserve"> Tablo g&#246;r&#252;nt&#252;leme,editleme</w:t></w:r>'

# --- Insert two new empty paragraphs between paragraph 3 (SE11/SE16N) and paragraph 4 (SE80) ---
$p3 = $d.Paragraphs.Item(3)
$ip = $d.Range($p3.Range.End - 1, $p3.Range.End - 1)
$ip.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item(3)
$ip = $d.Range($p3.Range.End - 1, $p3.Range.End - 1)
$ip.InsertParagraphAfter()

Write-Output "done"
